$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rename the "populationParamsFile" property to "populationsFile" and
# update its default value from "PopulationParameters.xlsx" to "Populations.xlsx"
$ws.Range("A6").Value = "populationsFile"
$ws.Range("B6").Value = "Populations.xlsx"

# Restore the active selection to reflect the edited cell
$ws.Range("B6").Select()
